$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: new entry "teste 2", inheriting the same wrap/vertical-top format
# as the rows above it (A2:A6) -- copy format from A6, then set the value.
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = "teste 2"

# Row 8: empty cell with an underlined font (new style).
$ws.Range("A8").Font.Underline = $true

# Page setup: portrait orientation (adds pageSetup element / printer info).
$ws.PageSetup.Orientation = 1

# Move the active selection down to the new last cell, A8.
$ws.Range("A8").Select() | Out-Null
